# Fruta / hortaliza, semanal
# A new week of data (previously at row 3, dated 2022-07-12 / serial 44754) is
# inserted as a new row, and the whole existing table (rows 3-22) shifts down
# one row (to rows 4-23). The new row 3 carries a new date/price observation
# (serial 44817) while every other row keeps its original data, just moved
# down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data rows (3-22) down by one row, opening up a fresh row 3.
$ws.Rows.Item(3).Insert()

# Populate the newly opened row 3 with the new weekly observation. All the
# "descriptive" columns are constant across every row in this table, so they
# are simply copied across; only the Fecha / Volumen / Precio columns differ.
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C3").Value = "Los Lagos"
$ws.Range("D3").Value = 44817
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 100112035
$ws.Range("G3").Value = "Bruselas (repollito)"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 90
$ws.Range("K3").Value = 23000
$ws.Range("L3").Value = 23000
$ws.Range("M3").Value = 23000
$ws.Range("N3").Value = "$/malla 15 kilos"
$ws.Range("O3").Value = "Provincia de Quillota"
$ws.Range("P3").Value = 1533
$ws.Range("Q3").Value = 15
$ws.Range("R3").Value = "Hortaliza"

# Give the new row the same date-format style as the rest of column D.
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat()
